# Update market-price derived columns (H:N) across leve-profit sheets
# to reflect refreshed pricing data from the scheduled market-board sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2755
$ws.Range("I70").Value = 2326.6667
$ws.Range("J70").Value = 2853.8462
$ws.Range("K70").Value = 6980.000100000001
$ws.Range("L70").Value = 8561.5386
$ws.Range("M70").Value = -6710.000100000001
$ws.Range("N70").Value = -9101.5386
$ws.Range("H73").Value = 2755
$ws.Range("I73").Value = 2326.6667
$ws.Range("J73").Value = 2853.8462
$ws.Range("K73").Value = 6980.000100000001
$ws.Range("L73").Value = 8561.5386
$ws.Range("M73").Value = -6044.000100000001
$ws.Range("N73").Value = -10433.5386
$ws.Range("H106").Value = 4180
$ws.Range("I106").Value = 3950
$ws.Range("J106").Value = 4333.3335
$ws.Range("K106").Value = 3950
$ws.Range("L106").Value = 4333.3335
$ws.Range("M106").Value = -3319
$ws.Range("N106").Value = -5595.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5299.956
$ws.Range("I32").Value = 3215.96
$ws.Range("K32").Value = 3215.96
$ws.Range("M32").Value = -2928.96
$ws.Range("H74").Value = 5004038
$ws.Range("I74").Value = 7408224.5
$ws.Range("J74").Value = 10727.23
$ws.Range("K74").Value = 7408224.5
$ws.Range("L74").Value = 10727.23
$ws.Range("M74").Value = -7407350.5
$ws.Range("N74").Value = -12475.23
$ws.Range("H77").Value = 5004038
$ws.Range("I77").Value = 7408224.5
$ws.Range("J77").Value = 10727.23
$ws.Range("K77").Value = 37041122.5
$ws.Range("L77").Value = 53636.14999999999
$ws.Range("M77").Value = -37036754.5
$ws.Range("N77").Value = -62372.14999999999
$ws.Range("H88").Value = 2086.8667
$ws.Range("I88").Value = 2238.125
$ws.Range("J88").Value = 1914
$ws.Range("K88").Value = 2238.125
$ws.Range("L88").Value = 1914
$ws.Range("M88").Value = -1832.125
$ws.Range("N88").Value = -2726
$ws.Range("H91").Value = 2086.8667
$ws.Range("I91").Value = 2238.125
$ws.Range("J91").Value = 1914
$ws.Range("K91").Value = 2238.125
$ws.Range("L91").Value = 1914
$ws.Range("M91").Value = -834.125
$ws.Range("N91").Value = -4722
$ws.Range("H110").Value = 3010.4
$ws.Range("I110").Value = 3642.2
$ws.Range("J110").Value = 2378.6
$ws.Range("K110").Value = 3642.2
$ws.Range("L110").Value = 2378.6
$ws.Range("M110").Value = -1597.2
$ws.Range("N110").Value = -6468.6
$ws.Range("H122").Value = 1145.7241
$ws.Range("I122").Value = 1157
$ws.Range("J122").Value = 1127.2727
$ws.Range("K122").Value = 3471
$ws.Range("L122").Value = 3381.8181
$ws.Range("M122").Value = -1021
$ws.Range("N122").Value = -8281.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 782.3333
$ws.Range("I25").Value = 753.5
$ws.Range("J25").Value = 840
$ws.Range("K25").Value = 753.5
$ws.Range("L25").Value = 840
$ws.Range("M25").Value = -518.5
$ws.Range("N25").Value = -1310

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1654.4412
$ws.Range("I31").Value = 1246.0588
$ws.Range("J31").Value = 2062.8235
$ws.Range("K31").Value = 1246.0588
$ws.Range("L31").Value = 2062.8235
$ws.Range("M31").Value = -951.0588
$ws.Range("N31").Value = -2652.8235
$ws.Range("H34").Value = 1654.4412
$ws.Range("I34").Value = 1246.0588
$ws.Range("J34").Value = 2062.8235
$ws.Range("K34").Value = 1246.0588
$ws.Range("L34").Value = 2062.8235
$ws.Range("M34").Value = -1044.0588
$ws.Range("N34").Value = -2466.8235
$ws.Range("H112").Value = 14900
$ws.Range("J112").Value = 14900
$ws.Range("L112").Value = 14900
$ws.Range("N112").Value = -17854
$ws.Range("H134").Value = 1889.3478
$ws.Range("I134").Value = 1996.9474
$ws.Range("J134").Value = 1378.25
$ws.Range("K134").Value = 5990.8422
$ws.Range("L134").Value = 4134.75
$ws.Range("M134").Value = -3455.8422
$ws.Range("N134").Value = -9204.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1865
$ws.Range("I97").Value = 500
$ws.Range("J97").Value = 2320
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 6960
$ws.Range("M97").Value = -1004
$ws.Range("N97").Value = -7952
$ws.Range("H98").Value = 649.4286
$ws.Range("I98").Value = 116
$ws.Range("J98").Value = 1049.5
$ws.Range("K98").Value = 348
$ws.Range("L98").Value = 3148.5
$ws.Range("M98").Value = 1150
$ws.Range("N98").Value = -6144.5
$ws.Range("H107").Value = 28463.31
$ws.Range("I107").Value = 66830.07000000001
$ws.Range("J107").Value = 18186.5
$ws.Range("K107").Value = 200490.21
$ws.Range("L107").Value = 54559.5
$ws.Range("M107").Value = -198570.21
$ws.Range("N107").Value = -58399.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21372
$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66864
$ws.Range("H70").Value = 4745.8
$ws.Range("I70").Value = 4695.4287
$ws.Range("J70").Value = 4863.3335
$ws.Range("K70").Value = 4695.4287
$ws.Range("L70").Value = 4863.3335
$ws.Range("M70").Value = -4425.4287
$ws.Range("N70").Value = -5403.3335
$ws.Range("H73").Value = 4745.8
$ws.Range("I73").Value = 4695.4287
$ws.Range("J73").Value = 4863.3335
$ws.Range("K73").Value = 4695.4287
$ws.Range("L73").Value = 4863.3335
$ws.Range("M73").Value = -3759.4287
$ws.Range("N73").Value = -6735.3335
$ws.Range("H80").Value = 113350.5
$ws.Range("I80").Value = 4375
$ws.Range("J80").Value = 186000.83
$ws.Range("K80").Value = 4375
$ws.Range("L80").Value = 186000.83
$ws.Range("M80").Value = -3377
$ws.Range("N80").Value = -187996.83
$ws.Range("H83").Value = 113350.5
$ws.Range("I83").Value = 4375
$ws.Range("J83").Value = 186000.83
$ws.Range("K83").Value = 21875
$ws.Range("L83").Value = 930004.1499999999
$ws.Range("M83").Value = -16883
$ws.Range("N83").Value = -939988.1499999999
$ws.Range("H113").Value = 4765.9062
$ws.Range("I113").Value = 6830.4736
$ws.Range("J113").Value = 1748.4615
$ws.Range("K113").Value = 6830.4736
$ws.Range("L113").Value = 1748.4615
$ws.Range("M113").Value = -4660.4736
$ws.Range("N113").Value = -6088.461499999999
$ws.Range("H122").Value = 3294.5
$ws.Range("I122").Value = 2290.48
$ws.Range("J122").Value = 4771
$ws.Range("K122").Value = 6871.440000000001
$ws.Range("L122").Value = 14313
$ws.Range("M122").Value = -4421.440000000001
$ws.Range("N122").Value = -19213

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 684.7692
$ws.Range("I46").Value = 711.1111
$ws.Range("J46").Value = 625.5
$ws.Range("K46").Value = 711.1111
$ws.Range("L46").Value = 625.5
$ws.Range("M46").Value = -523.1111
$ws.Range("N46").Value = -1001.5
$ws.Range("H61").Value = 1517.2142
$ws.Range("I61").Value = 1503.4166
$ws.Range("J61").Value = 1600
$ws.Range("K61").Value = 1503.4166
$ws.Range("L61").Value = 1600
$ws.Range("M61").Value = -1301.4166
$ws.Range("N61").Value = -2004
$ws.Range("H93").Value = 2709.3
$ws.Range("I93").Value = 2749.125
$ws.Range("K93").Value = 2749.125
$ws.Range("M93").Value = -1501.125
$ws.Range("H113").Value = 1517.2142
$ws.Range("I113").Value = 1503.4166
$ws.Range("J113").Value = 1600
$ws.Range("K113").Value = 1503.4166
$ws.Range("L113").Value = 1600
$ws.Range("M113").Value = 666.5834
$ws.Range("N113").Value = -5940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4159.5625
$ws.Range("I81").Value = 1076.625
$ws.Range("J81").Value = 7242.5
$ws.Range("K81").Value = 2153.25
$ws.Range("L81").Value = 14485
$ws.Range("M81").Value = -1092.25
$ws.Range("N81").Value = -16607
$ws.Range("H84").Value = 4159.5625
$ws.Range("I84").Value = 1076.625
$ws.Range("J84").Value = 7242.5
$ws.Range("K84").Value = 10766.25
$ws.Range("L84").Value = 72425
$ws.Range("M84").Value = -5462.25
$ws.Range("N84").Value = -83033
$ws.Range("H112").Value = 27852.428
$ws.Range("J112").Value = 27852.428
$ws.Range("L112").Value = 27852.428
$ws.Range("N112").Value = -30806.428
$ws.Range("H113").Value = 721.7
$ws.Range("I113").Value = 1062.25
$ws.Range("J113").Value = 332.5
$ws.Range("K113").Value = 3186.75
$ws.Range("L113").Value = 997.5
$ws.Range("M113").Value = -1016.75
$ws.Range("N113").Value = -5337.5
